$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: Title "A Table, with a caption"
# Split the merged "word + trailing space" runs back into separate
# word-run / space-run pairs (reverting the run-consolidation).
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Characters(1,1).Text = "A"
$title.Characters(2,1).Text = " "
$title.Characters(3,6).Text = "Table,"
$title.Characters(9,1).Text = " "
$title.Characters(10,4).Text = "with"
$title.Characters(14,1).Text = " "
$title.Characters(15,1).Text = "a"
$title.Characters(16,1).Text = " "

# Shape 3: Caption "Demonstration of simple table syntax, with alignment"
$caption = $s.Shapes.Item(3).TextFrame.TextRange
$caption.Characters(1,13).Text = "Demonstration"
$caption.Characters(14,1).Text = " "
$caption.Characters(15,2).Text = "of"
$caption.Characters(17,1).Text = " "
$caption.Characters(18,6).Text = "simple"
$caption.Characters(24,1).Text = " "
$caption.Characters(25,5).Text = "table"
$caption.Characters(30,1).Text = " "
$caption.Characters(31,7).Text = "syntax,"
$caption.Characters(38,1).Text = " "
$caption.Characters(39,4).Text = "with"
$caption.Characters(43,1).Text = " "
